$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.693.62"
$ws.Range("E2").Value = "  -4.55%  "

$ws.Range("D3").Value = "2.312.68"
$ws.Range("E3").Value = "  -6.01%  "

$c = $ws.Range("D4")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = $origStyle
$ws.Range("E4").Value = "  -0.03%  "

$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "305.85"
$c.Style = $origStyle
$ws.Range("E5").Value = "  -4.04%  "

$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "83.81"
$c.Style = $origStyle
$ws.Range("E6").Value = "  -8.20%  "

$c = $ws.Range("D7")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.528"
$c.Style = $origStyle
$ws.Range("E7").Value = "  -3.81%  "

$ws.Range("E8").Value = "  +0.05%  "

$c = $ws.Range("D9")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.479"
$c.Style = $origStyle
$ws.Range("E9").Value = "  -4.90%  "

$c = $ws.Range("D10")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0805"
$c.Style = $origStyle
$ws.Range("E10").Value = "  -5.41%  "

$c = $ws.Range("D11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "29.70"
$c.Style = $origStyle
$ws.Range("E11").Value = "  -9.03%  "

$ws.Range("E12").Value = "  +0.17%  "

$ws.Range("D13").Value = "2.674.64"
$ws.Range("E13").Value = "  -5.80%  "

$c = $ws.Range("D14")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.36"
$c.Style = $origStyle
$ws.Range("E14").Value = "  -7.10%  "

$c = $ws.Range("D15")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "14.54"
$c.Style = $origStyle
$ws.Range("E15").Value = "  -5.82%  "

$ws.Range("D16").Value = "2.322.33"
$ws.Range("E16").Value = "  -6.16%  "

$c = $ws.Range("D17")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.747"
$c.Style = $origStyle
$ws.Range("E17").Value = "  -4.84%  "

$ws.Range("D18").Value = "39.721.46"
$ws.Range("E18").Value = "  -4.34%  "

$ws.Range("D19").Value = "0.0₃0891"
$ws.Range("E19").Value = "  -4.87%  "

$c = $ws.Range("D20")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.04"
$c.Style = $origStyle
$ws.Range("E20").Value = "  -5.81%  "

$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "67.20"
$c.Style = $origStyle
$ws.Range("E21").Value = "  -6.76%  "

$c = $ws.Range("D22")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "10.51"
$c.Style = $origStyle
$ws.Range("E22").Value = "  -5.81%  "

$c = $ws.Range("D23")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "233.27"
$c.Style = $origStyle
$ws.Range("E23").Value = "  -1.99%  "

$c = $ws.Range("D24")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.53"
$c.Style = $origStyle
$ws.Range("E24").Value = "  -7.42%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("E26").Value = "  -7.28%  "

$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "22.93"
$c.Style = $origStyle
$ws.Range("E27").Value = "  -6.93%  "

$c = $ws.Range("D28")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.19"
$c.Style = $origStyle
$ws.Range("E28").Value = "  -1.97%  "

$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "9.18"
$c.Style = $origStyle
$ws.Range("E29").Value = "  -5.06%  "

$c = $ws.Range("D30")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "34.00"
$c.Style = $origStyle
$ws.Range("E30").Value = "  -5.74%  "

$c = $ws.Range("D31")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "149.71"
$c.Style = $origStyle
$ws.Range("E31").Value = "  -4.13%  "

$ws.Range("E32").Value = "  -0.11%  "

$c = $ws.Range("D33")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "5.05"
$c.Style = $origStyle
$ws.Range("E33").Value = "  -6.49%  "

$ws.Range("E34").Value = "  -4.82%  "

$c = $ws.Range("D35")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0714"
$c.Style = $origStyle
$ws.Range("E35").Value = "  -6.25%  "

$ws.Range("E36").Value = "  -2.50%  "

$c = $ws.Range("D37")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0988"
$c.Style = $origStyle
$ws.Range("E37").Value = "  -3.71%  "

$c = $ws.Range("D38")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.71"
$c.Style = $origStyle
$ws.Range("E38").Value = "  -6.54%  "

$c = $ws.Range("D39")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "15.40"
$c.Style = $origStyle
$ws.Range("E39").Value = "  -8.69%  "

$c = $ws.Range("D40")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.69"
$c.Style = $origStyle
$ws.Range("E40").Value = "  -7.22%  "

$c = $ws.Range("D41")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "3.76"
$c.Style = $origStyle
$ws.Range("E41").Value = "  -5.62%  "

$ws.Range("E42").Value = "  -2.52%  "

$ws.Range("D43").Value = "1.934.59"
$ws.Range("E43").Value = "  -3.27%  "

$c = $ws.Range("D44")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.0263"
$c.Style = $origStyle
$ws.Range("E44").Value = "  -6.28%  "

$c = $ws.Range("D45")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "17.52"
$c.Style = $origStyle
$ws.Range("E45").Value = "  -5.25%  "

$c = $ws.Range("D46")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "9.33"
$c.Style = $origStyle
$ws.Range("E46").Value = "  -1.97%  "

$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "2.65"
$c.Style = $origStyle
$ws.Range("E47").Value = "  -9.78%  "

$ws.Range("D48").Value = "2.542.07"
$ws.Range("E48").Value = "  -6.65%  "

$c = $ws.Range("D49")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "91.77"
$c.Style = $origStyle
$ws.Range("E49").Value = "  -5.17%  "

$c = $ws.Range("D50")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "69.72"
$c.Style = $origStyle
$ws.Range("E50").Value = "  -7.96%  "

$c = $ws.Range("D51")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "63.03"
$c.Style = $origStyle
$ws.Range("E51").Value = "  -5.18%  "
